$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new blood donation record as row 3
$ws.Range("A3").Value = "Kaviya"
$ws.Range("B3").Value = "A+"
$ws.Range("C3").Value = 19

# D3 holds a long numeric-looking id that must be stored as text (like D2 = "1"),
# so force text formatting while assigning it, then drop the formatting again so
# the cell ends up using the default (unstyled) cell format, same as D2.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1234567891"
$ws.Range("D3").ClearFormats()

# E3 holds the donation date (serial 44404 = 2021-07-27), formatted as yyyy-mm-dd
$ws.Range("E3").NumberFormat = "yyyy-mm-dd"
$ws.Range("E3").Value = 44404
